# Deploying to gh-pages from  @ 03fb76b9ac146585df05395f980ae353fb99e762
# Adds a "2020" column (Q) of data to the 1.3.1 SDG indicator table, mirroring
# the existing 2019 (P) column's formatting, and bolds the total row's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Header cells: clone formatting from column P (2019) into column Q (2020) ---
$headerRows = @(2, 3)
foreach ($r in $headerRows) {
    $src = $ws.Range("P" + $r)
    $dst = $ws.Range("Q" + $r)
    $src.Copy() | Out-Null
    $dst.PasteSpecial($xlPasteFormats) | Out-Null
}
$ws.Range("Q2").Borders.Item(9).Weight = -4138

# Year label for the new column
$ws.Range("Q3").Value = 2020

# --- Data rows 4-13: clone formatting from column P, then write the 2020 values ---
$dataRows = 4..13
foreach ($r in $dataRows) {
    $src = $ws.Range("P" + $r)
    $dst = $ws.Range("Q" + $r)
    $src.Copy() | Out-Null
    $dst.PasteSpecial($xlPasteFormats) | Out-Null
}

$ws.Range("Q4").Value = 17.7
$ws.Range("Q5").Value = 1.7006983633535606
$ws.Range("Q6").Value = 4.0792532187560786
$ws.Range("Q7").Value = 1.4965639329659175
$ws.Range("Q8").Value = 1.2345401844834025
$ws.Range("Q9").Value = 3.9182419607753913
$ws.Range("Q10").Value = 0.84723021008759791
$ws.Range("Q11").Value = 2.085763280904978
$ws.Range("Q12").Value = 1.8003095767645958
$ws.Range("Q13").Value = 0.49216467627561039

$excel.CutCopyMode = 0

# --- Row 4 is the combined "total" figure -- its numbers are bold in the new layout ---
$ws.Range("D4:Q4").Font.Bold = $true
